$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old bottom border on row 27 (previously the last "N" row)
$ws.Range("A27:B27").Borders.Item(9).LineStyle = -4142

# Set cell values for rows 2 through 33 (data rows)
$ws.Range("A2").Value = "group_1"
$ws.Range("B2").Value = 2.7041931202802711
$ws.Range("A3").Value = "group_2"
$ws.Range("B3").Value = 2.4982583762870894
$ws.Range("A4").Value = "difference"
$ws.Range("B4").Value = 0.20593474399318179
$ws.Range("A5").Value = "endowments"
$ws.Range("B5").Value = -0.023448545331923754
$ws.Range("A6").Value = "coefficients"
$ws.Range("B6").Value = 0.24455519635574152
$ws.Range("A7").Value = "interaction"
$ws.Range("B7").Value = -0.015171907030635978
$ws.Range("A8").Value = "age"
$ws.Range("B8").Value = 0.014453188373285667
$ws.Range("A9").Value = "LTHS"
$ws.Range("B9").Value = -0.025132750522614638
$ws.Range("A10").Value = "some_college"
$ws.Range("B10").Value = 0.0067974883864450868
$ws.Range("A11").Value = "college"
$ws.Range("B11").Value = -0.00066804655671171248
$ws.Range("A12").Value = "high_school"
$ws.Range("B12").Value = -0.005054431556861254
$ws.Range("A13").Value = "advanced_degree"
$ws.Range("B13").Value = -0.01011900567441181
$ws.Range("A14").Value = "foreign_born"
$ws.Range("B14").Value = -0.0018624938905276477
$ws.Range("A15").Value = "native"
$ws.Range("B15").Value = -0.0018624938905276488
$ws.Range("A16").Value = "age"
$ws.Range("B16").Value = -0.032229482695173325
$ws.Range("A17").Value = "LTHS"
$ws.Range("B17").Value = -0.015046868402275202
$ws.Range("A18").Value = "some_college"
$ws.Range("B18").Value = 0.0039833232039781924
$ws.Range("A19").Value = "college"
$ws.Range("B19").Value = -0.0068764829693986878
$ws.Range("A20").Value = "high_school"
$ws.Range("B20").Value = 0.0070835695284810945
$ws.Range("A21").Value = "advanced_degree"
$ws.Range("B21").Value = 0.0055836836672005707
$ws.Range("A22").Value = "foreign_born"
$ws.Range("B22").Value = -0.02240261832501322
$ws.Range("A23").Value = "native"
$ws.Range("B23").Value = 0.020610408859012164
$ws.Range("A24").Value = "Intercept"
$ws.Range("B24").Value = 0.28384966348893004
$ws.Range("A25").Value = "age"
$ws.Range("B25").Value = -0.0013237163090532855
$ws.Range("A26").Value = "LTHS"
$ws.Range("B26").Value = -0.0044595196330836943
$ws.Range("A27").Value = "some_college"
$ws.Range("B27").Value = -0.00094841028666147462
$ws.Range("A28").Value = "college"
$ws.Range("B28").Value = 0.0002401311513123358
$ws.Range("A29").Value = "high_school"
$ws.Range("B29").Value = 0.00046063150543630178
$ws.Range("A30").Value = "advanced_degree"
$ws.Range("B30").Value = -0.0018299467480741368
$ws.Range("A31").Value = "foreign_born"
$ws.Range("B31").Value = -0.0036555383552561234
$ws.Range("A32").Value = "native"
$ws.Range("B32").Value = -0.003655538355256126
$ws.Range("A33").Value = "N"
$ws.Range("B33").Value = 666

# Apply number formats and alignment for data rows 2-32 (coefficients) and row 33 (N, integer)
$ws.Range("A2:A33").HorizontalAlignment = -4131
$ws.Range("B2:B32").NumberFormat = "0.######"
$ws.Range("B2:B33").HorizontalAlignment = -4152
$ws.Range("B33").NumberFormat = "0"

# Bottom border under the final row (N)
$ws.Range("A33:B33").Borders.Item(9).LineStyle = 1
$ws.Range("A33:B33").Borders.Item(9).Weight = 2
